$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last tab (after the current last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Header formatting to match the other sheets' header style (bold, boxed, centered)
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").Borders.LineStyle = 1
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").VerticalAlignment = -4160

# Row 2 - MATCH_CODE 4727 (leading apostrophe keeps numeric-looking text as text,
# matching the sibling sheets where every exported field is a string).
# A lone apostrophe is Excel's text-quote escape for an empty string literal,
# matching the source export's blank fields.
$ws.Range("A2").Value = "'4727"
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NO"

# Row 3 - MATCH_CODE 4731
$ws.Range("A3").Value = "'4731"
$ws.Range("C3").Value = "'3"
$ws.Range("D3").Value = "'0"
$ws.Range("E3").Value = "'7.95%"
$ws.Range("F3").Value = "NO"

# BATTING_POSITION is a genuine number in the source data
$ws.Range("B3").Value = 2
